$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B (USGSCode) before the existing Latitude column,
# shifting Latitude/Longitude/Type/Name one column to the right.
$ws.Columns("B").Insert()
$ws.Range("B1").Value = "USGSCode"

# --- Fill in USGSCode values for stations that have a USGS gauge ---
$ws.Range("B2").Value = 11313440
$ws.Range("B5").Value = 11313434
$ws.Range("B8").Value = 11312676
$ws.Range("B12").Value = 11313433
$ws.Range("B14").Value = 11304810
$ws.Range("B16").Value = 11303500

# --- Update the Type column (now column E) for stations whose monitoring
#     type changed / was clarified ---
$ws.Range("E2").Value = "WQ, Flow, Velocity"
$ws.Range("E5").Value = "Flow"
$ws.Range("E8").Value = "WQ, Flow, Velocity"
$ws.Range("E12").Value = "WQ, Flow"
$ws.Range("E14").Value = "WQ, Flow"
$ws.Range("E16").Value = "Dayflow"

# --- Bold the StationCode for continuous (WQ-only, no gauge) stations ---
foreach ($addr in @("A2", "A3", "A4", "A6", "A7", "A13")) {
    $cell = $ws.Range($addr)
    $cell.Font.Bold = $true
    $cell.Font.Color = 0
    $cell.VerticalAlignment = -4108
}

# --- Add new Dayflow stations as additional rows ---
# First, copy the StationCode/USGSCode formatting (style used by the
# non-bold stations, e.g. row 5) down onto the new rows so the new cells
# pick up the same cell style instead of Excel minting a fresh one.
$ws.Range("A5:B5").Copy()
$ws.Range("A17:B19").PasteSpecial(-4122)

$ws.Range("A17").Value = "DTO"
$ws.Range("B17").Value = 11303500
$ws.Range("C17").Value = 38.059
$ws.Range("D17").Value = -122.025
$ws.Range("E17").Value = "Dayflow"
$ws.Range("F17").Value = "DELTA OUTFLOW"

$ws.Range("A18").Value = "HRO"
$ws.Range("C18").Value = 37.798
$ws.Range("D18").Value = -121.623
$ws.Range("E18").Value = "Dayflow"
$ws.Range("F18").Value = "HARVEY O BANKS PUMPING PLANT"

$ws.Range("A19").Value = "TRP"
$ws.Range("C19").Value = 37.8
$ws.Range("D19").Value = -121.585
$ws.Range("E19").Value = "Dayflow"
$ws.Range("F19").Value = "TRACY PUMPING PLANT"

# --- Remove one of the trailing placeholder rows (3 got consumed by new
#     data rows above, one fewer blank row remains at the bottom) ---
$ws.Rows("28").Delete()

# --- Column widths (bestFit-style widths on StationCode/USGSCode/Name) ---
# The engine quantizes ColumnWidth to 1/6-character steps, so these inputs
# are chosen to land on the stored-width value closest to the target.
$ws.Columns("A").ColumnWidth = 11
$ws.Columns("B").ColumnWidth = 11
$ws.Columns("E").ColumnWidth = 17.25

# --- Selection moved to G8 ---
$ws.Range("G8").Select()
